$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark that currently sits in paragraph 1 ---
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# --- 2. Paragraph 1: "Project Progress Report" -> center + underline ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.Alignment = 1
$p1.Range.Font.Underline = 1

# --- 3. Paragraph 2: "Hannah Weber, Jacob Paul, Marissa Kelley " -> center, trim trailing space ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.ParagraphFormat.Alignment = 1
$d.Content.Find.Execute("Hannah Weber, Jacob Paul, Marissa Kelley ", $true, $false, $false, $false, $false, $true, 1, $false, "Hannah Weber, Jacob Paul, Marissa Kelley", 2)

# --- 4. Insert new centered paragraph "12.7.18" after paragraph 2 ---
$p2 = $d.Paragraphs.Item(2)
$rng = $p2.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.InsertAfter("12.7.18")

Write-Host "Stage 1 complete"

# --- 5. Clear the old "We began our project..." paragraph (now paragraph 5), keep the (empty) paragraph ---
$p5 = $d.Paragraphs.Item(5)
$rng = $p5.Range
$rng.End = $rng.End - 1
$rng.Text = ""

Write-Host "Stage 2 complete"

# --- 6. Rewrite "We are planning on doing an oral presentation.  " paragraph (now paragraph 6)
#        into the new, longer "We began our project..." paragraph, with new formatting ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.ParagraphFormat.FirstLineIndent = 36
$p6.Range.ParagraphFormat.Alignment = 3

$rng = $p6.Range
$rng.End = $rng.End - 1
$newText = "We began our project by pulling the dataset from Kaggle. We found that it was quite large, so we began by doing a Sequel query to get the data out and also reduce the mass amount of data. We started by doing a groupby to look at which state had the most amount of wildfires, and turned it into an interactive bar charts. We’ve decided to predict looking at the causes of the wildfires, and look at it in the context of the state.  We haven’t run into any problems or challenges yet. Our next steps are to begin doing a random forest, perceptron and logistic regression on the data.  "
$rng.Text = $newText

Write-Host "Stage 3 complete"

# --- 7. Insert a new empty paragraph (justified) after the big paragraph ---
$p6 = $d.Paragraphs.Item(6)
$rng = $p6.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
# $rng is still collapsed right at the boundary (start of freshly-created empty paragraph 7)
$rng.ParagraphFormat.Alignment = 3

# move rng past the paragraph mark of the new empty paragraph 7, so the next
# InsertParagraphAfter/InsertAfter operate on paragraph 8 instead of paragraph 7
$rng.Start = $rng.Start + 1
$rng.End = $rng.End + 1

# --- 8. Insert the new "We are planning..." paragraph (justified) after that empty one ---
$rng.InsertParagraphAfter()
# $rng is now collapsed right at the boundary (start of freshly-created empty paragraph 8)
$rng.ParagraphFormat.Alignment = 3
$rng.InsertAfter("We are planning on doing an oral presentation, but if need be, we can just do a paper (if there are too many groups presenting). ")

# zero-length bookmark right after the text above, before the trailing two spaces
$bmRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

$rng.InsertAfter("  ")

Write-Host "Stage 4 complete"



